# Organizando las unidades desde las hojas de calculo
#
# Inserts a new "config" worksheet between "prop_mat" and "varios" that
# centralises scaling/unit variables previously scattered around the
# workbook, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

$propMat = $wb.Worksheets.Item("prop_mat")

# Create the new sheet right after "prop_mat" (i.e. right before "varios").
$config = $wb.Worksheets.Add($null, $propMat)
$config.Name = "config"

# --- Header row -----------------------------------------------------------
$config.Range("A1").Value = "variable"
$config.Range("B1").Value = "valor"

# Reuse the same bold/centred header look used on the other sheets
# (e.g. "prop_mat"!A1) instead of rebuilding it property-by-property, so we
# don't leave unused intermediate styles behind.
$propMat.Range("A1").Copy()
$config.Range("A1:B1").PasteSpecial(-4122)

# --- Data rows --------------------------------------------------------------
$config.Range("A2").Value = "esc_def"
$config.Range("B2").Value = 200

$config.Range("A3").Value = "esc_faxial"
$config.Range("B3").Value = 10

$config.Range("A4").Value = "esc_V"
$config.Range("B4").Value = 0.05

$config.Range("A5").Value = "esc_M"
$config.Range("B5").Value = 0.05

$config.Range("A6").Value = "U_FUER"
$config.Range("B6").Value = "kN"

$config.Range("A7").Value = "U_LONG"
$config.Range("B7").Value = "m"

# --- View state -------------------------------------------------------------
# "config" becomes the active sheet / active tab, taking the selection
# away from "prop_mat".
$config.Activate()
[void]$config.Range("D4").Select()
